$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9840301871299744
$ws.Range("B1").Value = 1.815277218818665
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.942726492881775
$ws.Range("E1").Value = 1.22331714630127
